# Quote's now have an itemized billing option.
# Update the "Itemized pricing for quotes" row's Status cell (C17) from
# "OPEN" to "Completed - Mar 16, 2010", and move the active selection to C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C17").Value = "Completed - Mar 16, 2010"

$ws.Range("C17").Select()
